$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of billing data
$ws.Range("A5").Value = (Get-Date -Year 2018 -Month 11 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "add plotly interaction"

# Move active selection to A6, as left by the editor after entering the row
$ws.Range("A6").Select()
